$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value is a plain decimal number look-alike
# (e.g. "223.07") must be forced to Text format first, mirroring the
# original inline-string/text storage, otherwise Excel auto-converts them
# to numeric values on assignment.
$ws.Range("D2").Value = "34.105.03"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").Value = "1.798.97"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.07"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -0.46%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.28"
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.285"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0715"
$ws.Range("E10").Value = "  +4.52%  "
$ws.Range("E11").Value = "  -1.14%  "
$ws.Range("D12").Value = "2.057.83"
$ws.Range("E12").Value = "  +0.65%  "
$ws.Range("D13").Value = "1.821.46"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.70"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.631"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("D16").Value = "34.131.29"
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.22"
$ws.Range("E17").Value = "  -1.27%  "
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "246.65"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.85"
$ws.Range("E22").Value = "  +4.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("E24").Value = "  -0.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.61"
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.57"
$ws.Range("E26").Value = "  +1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0525"
$ws.Range("E30").Value = "  +1.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.72"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +1.92%  "
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("D35").Value = "1.414.61"
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("E36").Value = "  +2.71%  "
$ws.Range("E37").Value = "  +0.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0187"
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.943"
$ws.Range("E39").Value = "  +5.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.24"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("E43").Value = "  +4.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.95"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0495"
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "107.34"
$ws.Range("E46").Value = "  +3.86%  "
$ws.Range("D47").Value = "1.956.84"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  -2.74%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.92"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").Value = "0.0₆0124"
$ws.Range("E51").Value = "  +2.46%  "
